# New weekly price record for "Coliflor" at Feria Lagunitas de Puerto Montt.
# A new row is inserted at row 229 (pushing the existing rows 229-302 down to
# 230-303), and the new row is populated with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(229).Insert()

$ws.Range("A229").Value = 4
$ws.Range("B229").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C229").Value = "Los Lagos"
$ws.Range("D229").Value = 44663
$ws.Range("E229").Value = 10
$ws.Range("F229").Value = 100112008
$ws.Range("G229").Value = "Coliflor"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 1000
$ws.Range("K229").Value = 1600
$ws.Range("L229").Value = 1700
$ws.Range("M229").Value = 1650
$ws.Range("N229").Value = "$/unidad"
$ws.Range("O229").Value = "Región Metropolitana"
$ws.Range("P229").Value = 1650
$ws.Range("Q229").Value = 1
$ws.Range("R229").Value = "Hortaliza"
